# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled-runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 889.1875
$ws.Range("I19").Value = 593.8333
$ws.Range("K19").Value = 593.8333
$ws.Range("M19").Value = -418.8333
$ws.Range("H42").Value = 1227.1666
$ws.Range("I42").Value = 35.285713
$ws.Range("K42").Value = 105.857139
$ws.Range("M42").Value = 124.142861
$ws.Range("H53").Value = 343.41177
$ws.Range("I53").Value = 277.22223
$ws.Range("J53").Value = 417.875
$ws.Range("K53").Value = 277.22223
$ws.Range("L53").Value = 417.875
$ws.Range("M53").Value = 359.77777
$ws.Range("N53").Value = -1691.875
$ws.Range("H92").Value = 2424.963
$ws.Range("I92").Value = 2073.85
$ws.Range("J92").Value = 3428.1428
$ws.Range("K92").Value = 2073.85
$ws.Range("L92").Value = 3428.1428
$ws.Range("M92").Value = -825.8499999999999
$ws.Range("N92").Value = -5924.1428
$ws.Range("H106").Value = 73336800
$ws.Range("I106").Value = 88002830
$ws.Range("K106").Value = 88002830
$ws.Range("M106").Value = -88002199
$ws.Range("H112").Value = 6601.75
$ws.Range("J112").Value = 7433.2
$ws.Range("L112").Value = 22299.6
$ws.Range("N112").Value = -24515.6
$ws.Range("H132").Value = 3124.587
$ws.Range("I132").Value = 2237.9302
$ws.Range("K132").Value = 6713.790599999999
$ws.Range("M132").Value = -4183.790599999999
$ws.Range("H135").Value = 511.82352
$ws.Range("I135").Value = 500.14285
$ws.Range("K135").Value = 4501.28565
$ws.Range("M135").Value = -1966.28565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3916.0588
$ws.Range("I122").Value = 3405
$ws.Range("K122").Value = 10215
$ws.Range("M122").Value = -7765
$ws.Range("H132").Value = 3007.3044
$ws.Range("I132").Value = 2325.818
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 6977.454000000001
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -4447.454000000001
$ws.Range("N132").Value = -59060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H33").Value = 65718.664
$ws.Range("I33").Value = 50000
$ws.Range("J33").Value = 68862.39999999999
$ws.Range("K33").Value = 50000
$ws.Range("L33").Value = 68862.39999999999
$ws.Range("M33").Value = -49664
$ws.Range("N33").Value = -69534.39999999999
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 23374.75
$ws.Range("J47").Value = 23374.75
$ws.Range("L47").Value = 23374.75
$ws.Range("N47").Value = -24506.75
$ws.Range("H58").Value = 2895.7368
$ws.Range("I58").Value = 3002.4443
$ws.Range("J58").Value = 2799.7
$ws.Range("K58").Value = 3002.4443
$ws.Range("L58").Value = 2799.7
$ws.Range("M58").Value = -2799.4443
$ws.Range("N58").Value = -3205.7
$ws.Range("H60").Value = 41527.832
$ws.Range("J60").Value = 46766.5
$ws.Range("L60").Value = 46766.5
$ws.Range("N60").Value = -47788.5
$ws.Range("H86").Value = 4417.391
$ws.Range("I86").Value = 4272.857
$ws.Range("J86").Value = 4642.222
$ws.Range("K86").Value = 4272.857
$ws.Range("L86").Value = 4642.222
$ws.Range("M86").Value = -3149.857
$ws.Range("N86").Value = -6888.222
$ws.Range("H89").Value = 4417.391
$ws.Range("I89").Value = 4272.857
$ws.Range("J89").Value = 4642.222
$ws.Range("K89").Value = 21364.285
$ws.Range("L89").Value = 23211.11
$ws.Range("M89").Value = -15748.285
$ws.Range("N89").Value = -34443.11
$ws.Range("H96").Value = 38570.5
$ws.Range("J96").Value = 38570.5
$ws.Range("L96").Value = 38570.5
$ws.Range("N96").Value = -44062.5
$ws.Range("H99").Value = 4732
$ws.Range("I99").Value = 4911.636
$ws.Range("J99").Value = 4485
$ws.Range("K99").Value = 4911.636
$ws.Range("L99").Value = 4485
$ws.Range("M99").Value = -3413.636
$ws.Range("N99").Value = -7481
$ws.Range("H122").Value = 3324.318
$ws.Range("I122").Value = 2698.25
$ws.Range("J122").Value = 4075.6
$ws.Range("K122").Value = 8094.75
$ws.Range("L122").Value = 12226.8
$ws.Range("M122").Value = -5644.75
$ws.Range("N122").Value = -17126.8
$ws.Range("H126").Value = 4732
$ws.Range("I126").Value = 4911.636
$ws.Range("J126").Value = 4485
$ws.Range("K126").Value = 14734.908
$ws.Range("L126").Value = 13455
$ws.Range("M126").Value = -12264.908
$ws.Range("N126").Value = -18395
$ws.Range("H136").Value = 2895.7368
$ws.Range("I136").Value = 3002.4443
$ws.Range("J136").Value = 2799.7
$ws.Range("K136").Value = 9007.332900000001
$ws.Range("L136").Value = 8399.099999999999
$ws.Range("M136").Value = -6457.332900000001
$ws.Range("N136").Value = -13499.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 628.56525
$ws.Range("I113").Value = 333.22223
$ws.Range("J113").Value = 818.4286
$ws.Range("K113").Value = 999.66669
$ws.Range("L113").Value = 2455.2858
$ws.Range("M113").Value = 1170.33331
$ws.Range("N113").Value = -6795.2858
$ws.Range("H121").Value = 1057323.6
$ws.Range("I121").Value = 7794.5
$ws.Range("J121").Value = 1267229.5
$ws.Range("K121").Value = 23383.5
$ws.Range("L121").Value = 3801688.5
$ws.Range("M121").Value = -22073.5
$ws.Range("N121").Value = -3804308.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1674.4524
$ws.Range("I102").Value = 1095.7812
$ws.Range("J102").Value = 3526.2
$ws.Range("K102").Value = 1095.7812
$ws.Range("L102").Value = 3526.2
$ws.Range("M102").Value = 526.2188000000001
$ws.Range("N102").Value = -6770.2
$ws.Range("H107").Value = 303
$ws.Range("I107").Value = 399.33334
$ws.Range("K107").Value = 399.33334
$ws.Range("M107").Value = 1520.66666
$ws.Range("H122").Value = 5768.5454
$ws.Range("I122").Value = 2858.3333
$ws.Range("K122").Value = 8574.999899999999
$ws.Range("M122").Value = -6124.999899999999
$ws.Range("H126").Value = 5785.7617
$ws.Range("I126").Value = 5719.75
$ws.Range("K126").Value = 17159.25
$ws.Range("M126").Value = -14689.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6809.2856
$ws.Range("I61").Value = 2213.5715
$ws.Range("J61").Value = 11405
$ws.Range("K61").Value = 2213.5715
$ws.Range("L61").Value = 11405
$ws.Range("M61").Value = -2011.5715
$ws.Range("N61").Value = -11809
$ws.Range("H113").Value = 6809.2856
$ws.Range("I113").Value = 2213.5715
$ws.Range("J113").Value = 11405
$ws.Range("K113").Value = 2213.5715
$ws.Range("L113").Value = 11405
$ws.Range("M113").Value = -43.57150000000001
$ws.Range("N113").Value = -15745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 75000
$ws.Range("J51").Value = 75000
$ws.Range("L51").Value = 75000
$ws.Range("N51").Value = -76020
$ws.Range("H53").Value = 49998.5
$ws.Range("J53").Value = 49998.5
$ws.Range("L53").Value = 49998.5
$ws.Range("N53").Value = -51212.5
$ws.Range("H62").Value = 6250
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 6333.3335
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 6333.3335
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -7581.3335
$ws.Range("H65").Value = 6250
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 6333.3335
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 31666.6675
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -37906.6675
$ws.Range("H81").Value = 1548
$ws.Range("I81").Value = 1876.4
$ws.Range("J81").Value = 1000.6667
$ws.Range("K81").Value = 3752.8
$ws.Range("L81").Value = 2001.3334
$ws.Range("M81").Value = -2691.8
$ws.Range("N81").Value = -4123.3334
$ws.Range("H84").Value = 1548
$ws.Range("I84").Value = 1876.4
$ws.Range("J84").Value = 1000.6667
$ws.Range("K84").Value = 18764
$ws.Range("L84").Value = 10006.667
$ws.Range("M84").Value = -13460
$ws.Range("N84").Value = -20614.667
